$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 19 (apartment "뉴타운1,2,3차", id 13922) -- all rows below shift up.
$ws.Rows.Item(19).Delete()

# Update the active selection to match the post-edit cursor position.
$ws.Range("F17").Select()
